# Updates the cryptocurrency price/volume table (D:E) to the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = "29.828.58"
$dCell.Style = "Normal"
$ws.Range("E2").Value = "  +2.51%  "

$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = "1.858.29"
$dCell.Style = "Normal"
$ws.Range("E3").Value = "  +2.00%  "

$dCell = $ws.Range("D4")
$dCell.NumberFormat = "@"
$dCell.Value = "0.9982"
$dCell.Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "

$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = "246.65"
$dCell.Style = "Normal"
$ws.Range("E5").Value = "  +2.09%  "

$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = "0.6390"
$dCell.Style = "Normal"
$ws.Range("E6").Value = "  +3.92%  "

$dCell = $ws.Range("D7")
$dCell.NumberFormat = "@"
$dCell.Value = "0.9988"
$dCell.Style = "Normal"
$ws.Range("E7").Value = "  -0.16%  "

$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"
$dCell.Value = "0.3010"
$dCell.Style = "Normal"
$ws.Range("E8").Value = "  +4.30%  "

$dCell = $ws.Range("D9")
$dCell.NumberFormat = "@"
$dCell.Value = "0.07520"
$dCell.Style = "Normal"
$ws.Range("E9").Value = "  +2.51%  "

$dCell = $ws.Range("D10")
$dCell.NumberFormat = "@"
$dCell.Value = "24.20"
$dCell.Style = "Normal"
$ws.Range("E10").Value = "  +5.24%  "

$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = "0.07681"
$dCell.Style = "Normal"
$ws.Range("E11").Value = "  +0.17%  "

$dCell = $ws.Range("D12")
$dCell.NumberFormat = "@"
$dCell.Value = "1.852.11"
$dCell.Style = "Normal"
$ws.Range("E12").Value = "  +2.17%  "

$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = "5.080"
$dCell.Style = "Normal"
$ws.Range("E13").Value = "  +2.69%  "

$dCell = $ws.Range("D14")
$dCell.NumberFormat = "@"
$dCell.Value = "0.6902"
$dCell.Style = "Normal"
$ws.Range("E14").Value = "  +4.13%  "

$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = "84.73"
$dCell.Style = "Normal"
$ws.Range("E15").Value = "  +3.96%  "

$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = "0.000009535"
$dCell.Style = "Normal"
$ws.Range("E16").Value = "  +6.41%  "

$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = "6.115"
$dCell.Style = "Normal"
$ws.Range("E17").Value = "  +4.05%  "

$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = "29.764.86"
$dCell.Style = "Normal"
$ws.Range("E18").Value = "  +2.37%  "

$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = "2.105.91"
$dCell.Style = "Normal"
$ws.Range("E19").Value = "  +3.23%  "

$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = "240.65"
$dCell.Style = "Normal"
$ws.Range("E20").Value = "  +1.81%  "

$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"
$dCell.Value = "12.68"
$dCell.Style = "Normal"
$ws.Range("E21").Value = "  +1.85%  "

$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = "0.9994"
$dCell.Style = "Normal"
$ws.Range("E22").Value = "  -0.09%  "

$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = "7.371"
$dCell.Style = "Normal"
$ws.Range("E23").Value = "  +3.45%  "

$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = "1.000"
$dCell.Style = "Normal"
$ws.Range("E24").Value = "  +0.07%  "

$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"
$dCell.Value = "159.92"
$dCell.Style = "Normal"
$ws.Range("E25").Value = "  +0.74%  "

$dCell = $ws.Range("D26")
$dCell.NumberFormat = "@"
$dCell.Value = "0.1425"
$dCell.Style = "Normal"
$ws.Range("E26").Value = "  +0.27%  "

$dCell = $ws.Range("D27")
$dCell.NumberFormat = "@"
$dCell.Value = "8.573"
$dCell.Style = "Normal"
$ws.Range("E27").Value = "  +1.69%  "

$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value = "17.98"
$dCell.Style = "Normal"
$ws.Range("E28").Value = "  +2.01%  "

$dCell = $ws.Range("D29")
$dCell.NumberFormat = "@"
$dCell.Value = "1.509"
$dCell.Style = "Normal"
$ws.Range("E29").Value = "  +1.90%  "

$dCell = $ws.Range("D30")
$dCell.NumberFormat = "@"
$dCell.Value = "0.06049"
$dCell.Style = "Normal"
$ws.Range("E30").Value = "  +8.13%  "

$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value = "1.262"
$dCell.Style = "Normal"
$ws.Range("E31").Value = "  +4.52%  "

$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = "4.154"
$dCell.Style = "Normal"
$ws.Range("E32").Value = "  +1.48%  "

$dCell = $ws.Range("D33")
$dCell.NumberFormat = "@"
$dCell.Value = "4.150"
$dCell.Style = "Normal"
$ws.Range("E33").Value = "  +1.24%  "

$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = "1.882"
$dCell.Style = "Normal"
$ws.Range("E34").Value = "  +2.98%  "

$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"
$dCell.Value = "1.161"
$dCell.Style = "Normal"
$ws.Range("E35").Value = "  +2.58%  "

$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = "0.7360"
$dCell.Style = "Normal"
$ws.Range("E36").Value = "  -0.01%  "

$dCell = $ws.Range("D37")
$dCell.NumberFormat = "@"
$dCell.Value = "2.617"
$dCell.Style = "Normal"
$ws.Range("E37").Value = "  -0.31%  "

$ws.Range("E38").Value = "  +1.59%  "

$dCell = $ws.Range("D39")
$dCell.NumberFormat = "@"
$dCell.Value = "1.228.53"
$dCell.Style = "Normal"
$ws.Range("E39").Value = "  +2.11%  "

$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = "0.01793"
$dCell.Style = "Normal"
$ws.Range("E40").Value = "  +2.07%  "

$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = "6.395"
$dCell.Style = "Normal"
$ws.Range("E41").Value = "  +0.44%  "

$dCell = $ws.Range("D42")
$dCell.NumberFormat = "@"
$dCell.Value = "0.9232"
$dCell.Style = "Normal"
$ws.Range("E42").Value = "  +3.56%  "

$ws.Range("E43").Value = "  +0.09%  "

$dCell = $ws.Range("D44")
$dCell.NumberFormat = "@"
$dCell.Value = "2.017.99"
$dCell.Style = "Normal"
$ws.Range("E44").Value = "  +3.74%  "

$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = "102.56"
$dCell.Style = "Normal"
$ws.Range("E45").Value = "  +1.76%  "

$dCell = $ws.Range("D46")
$dCell.NumberFormat = "@"
$dCell.Value = "66.60"
$dCell.Style = "Normal"
$ws.Range("E46").Value = "  +3.09%  "

$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value = "0.00000000123"
$dCell.Style = "Normal"
$ws.Range("E47").Value = "  +2.58%  "

$dCell = $ws.Range("D48")
$dCell.NumberFormat = "@"
$dCell.Value = "0.5085"
$dCell.Style = "Normal"
$ws.Range("E48").Value = "  +0.22%  "

$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"
$dCell.Value = "9.367"
$dCell.Style = "Normal"
$ws.Range("E49").Value = "  +3.03%  "

$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value = "0.4098"
$dCell.Style = "Normal"
$ws.Range("E50").Value = "  +2.60%  "

$dCell = $ws.Range("D51")
$dCell.NumberFormat = "@"
$dCell.Value = "0.1150"
$dCell.Style = "Normal"
$ws.Range("E51").Value = "  +3.22%  "
